$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace header row with single title cell, clear B1:L1
$ws.Range("A1").Value = "Herald College Kathmandu"
$ws.Range("B1:L1").ClearContents()

# Data rows 2-10: remap columns (Hours/Code/Title shift, Room/Group swap) and drop K,L
# Row 2
$ws.Range("A2").Value = "SUN"
$ws.Range("B2").Value = "7:00-9:30"
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = "5CS022"
$ws.Range("E2").Value = "Human Computer Interaction"
$ws.Range("F2").Value = "Workshop"
$ws.Range("G2").Value = "Mr. Ayush Shakya"
$ws.Range("H2").Value = "L5CG12"
$ws.Range("I2").Value = "WLV"
$ws.Range("J2").Value = "Lab-01 Mander"
$ws.Range("K2:L2").ClearContents()

# Row 3
$ws.Range("A3").Value = "SUN"
$ws.Range("B3").Value = "10:00-12:00"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "5CS020"
$ws.Range("E3").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F3").Value = "Tutorial"
$ws.Range("G3").Value = "Mr. Prabin Sapkota"
$ws.Range("H3").Value = "L5CG12"
$ws.Range("I3").Value = "HCK"
$ws.Range("J3").Value = "TR-11 Nagarjung"
$ws.Range("K3:L3").ClearContents()

# Row 4
$ws.Range("A4").Value = "MON"
$ws.Range("B4").Value = "13:00-15:30"
$ws.Range("C4").Value = 2.5
$ws.Range("D4").Value = "5CS020"
$ws.Range("E4").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F4").Value = "Workshop"
$ws.Range("G4").Value = "Mr. Prabin Sapkota"
$ws.Range("H4").Value = "L5CG12"
$ws.Range("I4").Value = "WLV"
$ws.Range("J4").Value = "Lab-01 Mander"
$ws.Range("K4:L4").ClearContents()

# Row 5
$ws.Range("A5").Value = "TUE"
$ws.Range("B5").Value = "7:00-9:00"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "5CS024"
$ws.Range("E5").Value = "Collaborative Development"
$ws.Range("F5").Value = "Lecture"
$ws.Range("G5").Value = "Mr. Raj Shrestha"
$ws.Range("H5").Value = "L5CG(12+13+14)"
$ws.Range("I5").Value = "WLV"
$ws.Range("J5").Value = "LT-03 Walsall"
$ws.Range("K5:L5").ClearContents()

# Row 6
$ws.Range("A6").Value = "WED"
$ws.Range("B6").Value = "7:00-9:00"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "5CS022"
$ws.Range("E6").Value = "Human Computer Interaction"
$ws.Range("F6").Value = "Lecture"
$ws.Range("G6").Value = "Mr. Ayush Shakya"
$ws.Range("H6").Value = "L5CG(12+13+14)"
$ws.Range("I6").Value = "WLV"
$ws.Range("J6").Value = "LT-01 Wulfruna"
$ws.Range("K6:L6").ClearContents()

# Row 7
$ws.Range("A7").Value = "WED"
$ws.Range("B7").Value = "9:30-11:30"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "5CS024"
$ws.Range("E7").Value = "Collaborative Development"
$ws.Range("F7").Value = "Tutorial"
$ws.Range("G7").Value = "Mr. Anmol Adhikari"
$ws.Range("H7").Value = "L5CG12"
$ws.Range("I7").Value = "WLV"
$ws.Range("J7").Value = "SR-04 Crompton"
$ws.Range("K7:L7").ClearContents()

# Row 8
$ws.Range("A8").Value = "THU"
$ws.Range("B8").Value = "9:30-11:30"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "5CS020"
$ws.Range("E8").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F8").Value = "Lecture"
$ws.Range("G8").Value = "Mr. Sumanta Silwal"
$ws.Range("H8").Value = "L5CG(12+13+14)"
$ws.Range("I8").Value = "WLV"
$ws.Range("J8").Value = "LT-01 Wulfruna"
$ws.Range("K8:L8").ClearContents()

# Row 9
$ws.Range("A9").Value = "THU"
$ws.Range("B9").Value = "12:30-15:00"
$ws.Range("C9").Value = 2.5
$ws.Range("D9").Value = "5CS024"
$ws.Range("E9").Value = "Collaborative Development"
$ws.Range("F9").Value = "Workshop"
$ws.Range("G9").Value = "Mr. Anmol Adhikari"
$ws.Range("H9").Value = "L5CG12"
$ws.Range("I9").Value = "WLV"
$ws.Range("J9").Value = "SR-03 Wolves"
$ws.Range("K9:L9").ClearContents()

# Row 10
$ws.Range("A10").Value = "FRI"
$ws.Range("B10").Value = "7:00-9:00"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "5CS022"
$ws.Range("E10").Value = "Human Computer Interaction"
$ws.Range("F10").Value = "Tutorial"
$ws.Range("G10").Value = "Mr. Ayush Shakya"
$ws.Range("H10").Value = "L5CG12"
$ws.Range("I10").Value = "WLV"
$ws.Range("J10").Value = "SR-03 Wolves"
$ws.Range("K10:L10").ClearContents()

